$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17, 18 and 19 hold three "artfynd" observation records that need to be
# cyclically rotated: row 17 must end up with the data that currently lives
# in row 18, row 18 must end up with the data currently in row 19, and row 19
# must end up with the data currently in row 17 (A17/A18/A19 ids
# 131244253/131244274/131244273 become 131244274/131244273/131244253).
#
# Only the columns below actually differ between the three rows; every other
# column (D, I, K, N, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY, ...)
# already holds the same value in all three rows, so it is left untouched.
$cols = @("A","B","E","F","G","H","J","L","M","Q","R","AC","AF","AJ","AK","AO")

# Snapshot the current ("before") values for the three rows, column by column.
$before17 = @{}
$before18 = @{}
$before19 = @{}
foreach ($col in $cols) {
    $before17[$col] = $ws.Range("$col" + "17").Value2
    $before18[$col] = $ws.Range("$col" + "18").Value2
    $before19[$col] = $ws.Range("$col" + "19").Value2
}

# Write the rotated values: 17 <- 18, 18 <- 19, 19 <- 17 (original values).
foreach ($col in $cols) {
    $ws.Range("$col" + "17").Value2 = $before18[$col]
    $ws.Range("$col" + "18").Value2 = $before19[$col]
    $ws.Range("$col" + "19").Value2 = $before17[$col]
}
